$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19
$ws.Range("H19").Value = 1673304
$ws.Range("I19").Value = 1617.9166
$ws.Range("J19").Value = 3496961.5
$ws.Range("K19").Value = 1617.9166
$ws.Range("L19").Value = 3496961.5
$ws.Range("M19").Value = -1442.9166
$ws.Range("N19").Value = -3497311.5

# ALC row 86
$ws.Range("H86").Value = 2520.75
$ws.Range("I86").Value = 2150
$ws.Range("J86").Value = 3039.8
$ws.Range("K86").Value = 2150
$ws.Range("L86").Value = 3039.8
$ws.Range("M86").Value = -1027
$ws.Range("N86").Value = -5285.8

# ALC row 89
$ws.Range("H89").Value = 2520.75
$ws.Range("I89").Value = 2150
$ws.Range("J89").Value = 3039.8
$ws.Range("K89").Value = 10750
$ws.Range("L89").Value = 15199
$ws.Range("M89").Value = -5134
$ws.Range("N89").Value = -26431

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 2315.7778
$ws.Range("I2").Value = 2453.3684
$ws.Range("K2").Value = 2453.3684
$ws.Range("M2").Value = -2340.3684

# ARM row 61
$ws.Range("H61").Value = 1970.091
$ws.Range("I61").Value = 1611.5
$ws.Range("J61").Value = 3302
$ws.Range("K61").Value = 1611.5
$ws.Range("L61").Value = 3302
$ws.Range("M61").Value = -1399.5
$ws.Range("N61").Value = -3726

# ARM row 97
$ws.Range("H97").Value = 46776.59
$ws.Range("I97").Value = 63629.625
$ws.Range("J97").Value = 1835.1666
$ws.Range("K97").Value = 63629.625
$ws.Range("L97").Value = 1835.1666
$ws.Range("M97").Value = -63133.625
$ws.Range("N97").Value = -2827.1666

# ARM row 116
$ws.Range("H116").Value = 2315.7778
$ws.Range("I116").Value = 2453.3684
$ws.Range("K116").Value = 2453.3684
$ws.Range("M116").Value = -159.3683999999998

# ARM row 132
$ws.Range("H132").Value = 2190.111
$ws.Range("I132").Value = 1263.2916
$ws.Range("J132").Value = 4043.75
$ws.Range("K132").Value = 3789.8748
$ws.Range("L132").Value = 12131.25
$ws.Range("M132").Value = -1259.8748
$ws.Range("N132").Value = -17191.25

# ARM row 136
$ws.Range("H136").Value = 1970.091
$ws.Range("I136").Value = 1611.5
$ws.Range("J136").Value = 3302
$ws.Range("K136").Value = 4834.5
$ws.Range("L136").Value = 9906
$ws.Range("M136").Value = -2284.5
$ws.Range("N136").Value = -15006

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 2315.7778
$ws.Range("I3").Value = 2453.3684
$ws.Range("K3").Value = 2453.3684
$ws.Range("M3").Value = -2339.3684

# BSM row 86
$ws.Range("H86").Value = 2027.6666
$ws.Range("I86").Value = 1740.8948
$ws.Range("J86").Value = 2708.75
$ws.Range("K86").Value = 1740.8948
$ws.Range("L86").Value = 2708.75
$ws.Range("M86").Value = -617.8948
$ws.Range("N86").Value = -4954.75

# BSM row 89
$ws.Range("H89").Value = 2027.6666
$ws.Range("I89").Value = 1740.8948
$ws.Range("J89").Value = 2708.75
$ws.Range("K89").Value = 8704.474
$ws.Range("L89").Value = 13543.75
$ws.Range("M89").Value = -3088.474
$ws.Range("N89").Value = -24775.75

# BSM row 107
$ws.Range("H107").Value = 1005.8
$ws.Range("I107").Value = 1070.091
$ws.Range("J107").Value = 829
$ws.Range("K107").Value = 1070.091
$ws.Range("L107").Value = 829
$ws.Range("M107").Value = 849.9090000000001
$ws.Range("N107").Value = -4669

# BSM row 134
$ws.Range("H134").Value = 5001636
$ws.Range("I134").Value = 7813838.5
$ws.Range("J134").Value = 2165.4443
$ws.Range("K134").Value = 23441515.5
$ws.Range("L134").Value = 6496.3329
$ws.Range("M134").Value = -23438980.5
$ws.Range("N134").Value = -11566.3329

# BSM row 135
$ws.Range("H135").Value = 44499
$ws.Range("J135").Value = 47843.547
$ws.Range("L135").Value = 47843.547
$ws.Range("N135").Value = -57983.547

$ws = $wb.Worksheets.Item("CRP")
# CRP row 86
$ws.Range("H86").Value = 25964.38
$ws.Range("I86").Value = 39945.77
$ws.Range("J86").Value = 3244.625
$ws.Range("K86").Value = 39945.77
$ws.Range("L86").Value = 3244.625
$ws.Range("M86").Value = -38822.77
$ws.Range("N86").Value = -5490.625

# CRP row 89
$ws.Range("H89").Value = 25964.38
$ws.Range("I89").Value = 39945.77
$ws.Range("J89").Value = 3244.625
$ws.Range("K89").Value = 199728.85
$ws.Range("L89").Value = 16223.125
$ws.Range("M89").Value = -194112.85
$ws.Range("N89").Value = -27455.125

# CRP row 134
$ws.Range("H134").Value = 18520198
$ws.Range("I134").Value = 35716000
$ws.Range("J134").Value = 1642
$ws.Range("K134").Value = 107148000
$ws.Range("L134").Value = 4926
$ws.Range("M134").Value = -107145465
$ws.Range("N134").Value = -9996

$ws = $wb.Worksheets.Item("CUL")
# CUL row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

# CUL row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM row 40
$ws.Range("H40").Value = 14259
$ws.Range("J40").Value = 14259
$ws.Range("L40").Value = 14259
$ws.Range("N40").Value = -14561

# GSM row 43
$ws.Range("H43").Value = 7260.091
$ws.Range("I43").Value = 3983
$ws.Range("J43").Value = 12995
$ws.Range("K43").Value = 3983
$ws.Range("L43").Value = 12995
$ws.Range("M43").Value = -3832
$ws.Range("N43").Value = -13297

# GSM row 132
$ws.Range("H132").Value = 18526492
$ws.Range("I132").Value = 30313222
$ws.Range("J132").Value = 4489.2856
$ws.Range("K132").Value = 90939666
$ws.Range("L132").Value = 13467.8568
$ws.Range("M132").Value = -90937136
$ws.Range("N132").Value = -18527.8568

$ws = $wb.Worksheets.Item("LTW")
# LTW row 18
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1828

# LTW row 22
$ws.Range("H22").Value = 652.7273
$ws.Range("I22").Value = 335
$ws.Range("J22").Value = 723.3333
$ws.Range("K22").Value = 335
$ws.Range("L22").Value = 723.3333
$ws.Range("M22").Value = -40
$ws.Range("N22").Value = -1313.3333

# LTW row 27
$ws.Range("H27").Value = 652.7273
$ws.Range("I27").Value = 335
$ws.Range("J27").Value = 723.3333
$ws.Range("K27").Value = 335
$ws.Range("L27").Value = 723.3333
$ws.Range("M27").Value = -228
$ws.Range("N27").Value = -937.3333

# LTW row 46
$ws.Range("H46").Value = 853.34485
$ws.Range("I46").Value = 807.35
$ws.Range("J46").Value = 955.55554
$ws.Range("K46").Value = 807.35
$ws.Range("L46").Value = 955.55554
$ws.Range("M46").Value = -619.35
$ws.Range("N46").Value = -1331.55554

# LTW row 61
$ws.Range("H61").Value = 3846.875
$ws.Range("I61").Value = 4222
$ws.Range("J61").Value = 3793.2856
$ws.Range("K61").Value = 4222
$ws.Range("L61").Value = 3793.2856
$ws.Range("M61").Value = -4020
$ws.Range("N61").Value = -4197.2856

# LTW row 113
$ws.Range("H113").Value = 3846.875
$ws.Range("I113").Value = 4222
$ws.Range("J113").Value = 3793.2856
$ws.Range("K113").Value = 4222
$ws.Range("L113").Value = 3793.2856
$ws.Range("M113").Value = -2052
$ws.Range("N113").Value = -8133.2856

$ws = $wb.Worksheets.Item("WVR")
# WVR row 100
$ws.Range("H100").Value = 5682780
$ws.Range("I100").Value = 9091808
$ws.Range("J100").Value = 1066.6666
$ws.Range("K100").Value = 18183616
$ws.Range("L100").Value = 2133.3332
$ws.Range("M100").Value = -18183075
$ws.Range("N100").Value = -3215.3332

# WVR row 122
$ws.Range("H122").Value = 1951.5555
$ws.Range("I122").Value = 1427.3334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4282.0002
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1832.0002
$ws.Range("N122").Value = -13900

# WVR row 126
$ws.Range("H126").Value = 1977.6818
$ws.Range("I126").Value = 1489.5625
$ws.Range("J126").Value = 3279.3333
$ws.Range("K126").Value = 4468.6875
$ws.Range("L126").Value = 9837.999899999999
$ws.Range("M126").Value = -1998.6875
$ws.Range("N126").Value = -14777.9999
